$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 553, shifting existing rows 553:609 down to 554:610
$ws.Rows.Item(553).Insert()

# Populate the newly inserted row 553 with its values
$ws.Cells.Item(553, 1).Value = 4
$ws.Cells.Item(553, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(553, 3).Value = "Los Lagos"
$ws.Cells.Item(553, 4).Value = (Get-Date -Year 2023 -Month 10 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(553, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(553, 5).Value = 10
$ws.Cells.Item(553, 6).Value = 100112023
$ws.Cells.Item(553, 7).Value = "Brócoli"
$ws.Cells.Item(553, 8).Value = "Sin especificar"
$ws.Cells.Item(553, 9).Value = "Primera"
$ws.Cells.Item(553, 10).Value = 1500
$ws.Cells.Item(553, 11).Value = 1500
$ws.Cells.Item(553, 12).Value = 1500
$ws.Cells.Item(553, 13).Value = 1500
$ws.Cells.Item(553, 14).Value = "`$/unidad"
$ws.Cells.Item(553, 15).Value = "Región Metropolitana"
$ws.Cells.Item(553, 16).Value = 1500
$ws.Cells.Item(553, 17).Value = 1
$ws.Cells.Item(553, 18).Value = "Hortaliza"
